$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text cells (safe from numeric auto-conversion): assign directly.
$ws.Range("D2").Value = "27.623.61"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").Value = "1.847.81"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  -0.45%  "
$ws.Range("E9").Value = "  +2.47%  "
$ws.Range("E10").Value = "  +1.03%  "
$ws.Range("E11").Value = "  -2.66%  "
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").Value = "1.885.78"
$ws.Range("E13").Value = "  +1.91%  "
$ws.Range("E14").Value = "  -0.91%  "
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("E16").Value = "  +0.74%  "
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("E18").Value = "  +2.79%  "
$ws.Range("E19").Value = "  +1.17%  "
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("E21").Value = "  -0.78%  "
$ws.Range("D22").Value = "27.643.88"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("E24").Value = "  -2.60%  "
$ws.Range("D25").Value = "2.064.22"
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("E26").Value = "  -2.96%  "
$ws.Range("E27").Value = "  +1.05%  "
$ws.Range("E28").Value = "  +2.41%  "
$ws.Range("E29").Value = "  +9.05%  "
$ws.Range("E30").Value = "  -1.20%  "
$ws.Range("E31").Value = "  +2.97%  "
$ws.Range("E32").Value = "  -0.36%  "
$ws.Range("E33").Value = "  -2.76%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("E34").Value = "  +2.56%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("E35").Value = "  +0.66%  "
$ws.Range("E36").Value = "  +0.94%  "
$ws.Range("E37").Value = "  -0.62%  "
$ws.Range("E38").Value = "  -0.91%  "
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("E40").Value = "  -4.91%  "
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("E42").Value = "  +0.74%  "
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("E44").Value = "  +0.61%  "
$ws.Range("E45").Value = "  -1.27%  "
$ws.Range("E46").Value = "  -0.73%  "
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("E50").Value = "  -1.52%  "
$ws.Range("E51").Value = "  -0.45%  "

# Numeric-looking text cells: force text format so Excel keeps the literal string
# (e.g. "1.000" does not become the number 1), then restore the default style
# so no stray style index is left attached to the cell.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4282"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.96"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07314"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.67"
$ws.Range("D12").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.518"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06912"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "79.92"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000009042"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.956"
$ws.Range("D23").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.989"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "121.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.283"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.847"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08914"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7643"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.997"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.553"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.102"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05413"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01935"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.818"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5077"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.754"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.365"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.06549"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.32"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "105.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4670"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.9999"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.622"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.54"
$ws.Range("D51").Style = "Normal"
